$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in newly recorded "Story Points (Done that day)" values for
# 08. and 09. February (column D and E) with 0, matching 07. Feb (C2) which
# already had a value.
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# Re-enter the "Ideal (Done that day)" row formula across the whole range in
# one shot, so Excel collapses it back into a single shared formula group.
$ws.Range("C4:W4").Formula = "=`$B`$3/COUNTA(`$C`$1:`$W`$1)"

# Move the active selection to E3 (as last edited/selected cell).
$ws.Range("E3").Select()
